# Update "想去人数" (F) and "最低票价" (G) figures on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 333
    $ws.Range("G2").Value = 65

    $ws.Range("F4").Value = 1508

    $ws.Range("F5").Value = 18

    $ws.Range("F6").Value = 42

    $ws.Range("F7").Value = 123

    $ws.Range("F8").Value = 49

    $ws.Range("F9").Value = 319
}
